$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B7").Value = 3647541.63
$ws.Range("C7").Value = -17.90507545701222
$ws.Range("D7").Value = 3203
$ws.Range("E7").Value = 3203
$ws.Range("F7").Value = 1138.789144551982
$ws.Range("G7").Value = 21.38668830333752
